# [Mod] Excel 수치 조정
# Adjust monster stat numbers on the "Monsters" sheet and switch the
# active/selected sheet from "Characters" to "Monsters".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monsters")

# Rows 2-10: def (col C) 10 -> 0, atkSpeed (col E) 3 -> 2
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 5).Value = 2
}

# Row 2 only: exp (col I) 3 -> 10
$ws.Cells.Item(2, 9).Value = 10

# Make "Monsters" the active/selected sheet with H12 selected
# (previously "Characters" was the selected tab with K10 selected there).
$ws.Activate()
$ws.Range("H12").Select()
